$d = $word.ActiveDocument

# The document's body paragraphs were reshuffled: a chain of paragraph
# "value" texts rotates into the next slot (the old Objetivos blurb moves
# to "Docente(s) Responsável(eis)", the old "A definir..." line moves up
# into Objetivos, etc.) while paragraph styles and any fixed label runs
# (Método:/Critério:/Norma de recuperação:) stay exactly where they are.
# We therefore just overwrite each paragraph's text content directly
# (by paragraph index) so formatting/styles are preserved, instead of
# trying to move paragraphs around.

# Paragraph 6: Objetivos (PT) body <- old "Programa resumido" (PT) text
$d.Paragraphs.Item(6).Range.Text = "A definir de acordo com o tópico programado"

# Paragraph 7: Objetivos (EN, italic) body <- old "Programa resumido" (EN) text
$d.Paragraphs.Item(7).Range.Text = "To be defined according to the scheduled topic"

# Paragraph 9: Docente(s) Responsável(eis) list value <- old Objetivos (PT) text
$d.Paragraphs.Item(9).Range.Text = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."

# Paragraph 11: Programa resumido (PT) body <- old "Programa" (PT) text
$d.Paragraphs.Item(11).Range.Text = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."

# Paragraph 12: Programa resumido (EN, italic) body <- old Objetivos (EN) text
$d.Paragraphs.Item(12).Range.Text = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art"

# Paragraph 14: Programa (PT) body <- old Avaliação "Método:" value text
$d.Paragraphs.Item(14).Range.Text = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."

# Paragraph 17 ("Avaliação" list): three labelled values (Método:/Critério:/
# Norma de recuperação:) each rotate to the next one. Locate the label runs
# by anchoring on their fixed label text (never changes), then overwrite
# just the value range that follows each label - right to left, so the
# not-yet-processed offsets earlier in the paragraph stay valid.
$p17 = $d.Paragraphs.Item(17)
$paraStart = $p17.Range.Start
$paraEndExclMark = $p17.Range.End - 1  # exclude the paragraph mark

$rMetodoLabel = $d.Range($paraStart, $paraEndExclMark)
$rMetodoLabel.Find.Execute("Método: ") | Out-Null

$rCriterioLabel = $d.Range($paraStart, $paraEndExclMark)
$rCriterioLabel.Find.Execute("Critério: ") | Out-Null

$rNormaLabel = $d.Range($paraStart, $paraEndExclMark)
$rNormaLabel.Find.Execute("Norma de recuperação: ") | Out-Null

$valueMetodoStart = $rMetodoLabel.End
$valueMetodoEnd = $rCriterioLabel.Start
$valueCriterioStart = $rCriterioLabel.End
$valueCriterioEnd = $rNormaLabel.Start
$valueNormaStart = $rNormaLabel.End
$valueNormaEnd = $paraEndExclMark

$brk = [char]11   # represents the <w:br/> line break inside a Range.Text value

# Right to left: Norma value, then Critério value, then Método value.
$rValueNorma = $d.Range($valueNormaStart, $valueNormaEnd)
$rValueNorma.Text = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."

$rValueCriterio = $d.Range($valueCriterioStart, $valueCriterioEnd)
$rValueCriterio.Text = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2" + $brk

$rValueMetodo = $d.Range($valueMetodoStart, $valueMetodoEnd)
$rValueMetodo.Text = "Média ponderada das avaliações (M)." + $brk

# Paragraph 19: Bibliografia body <- old "Docente(s)" list value
$d.Paragraphs.Item(19).Range.Text = "11079086 - Herlandí de Souza Andrade"
